$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly_mean (D) and weekly_share (F) values per row, 2020-10-23 update
$updates = @{
    2  = @{ D = 23.467;  F = 0.757 }
    3  = @{ D = 23.467;  F = 0.717 }
    4  = @{ D = 23.467;  F = 1.088 }
    5  = @{ D = 23.467;  F = 1.368 }
    6  = @{ D = 23.467;  F = 1.07 }
    7  = @{ D = 22.119;  F = 0.928 }
    8  = @{ D = 22.119;  F = 0.852 }
    9  = @{ D = 22.119;  F = 1.305 }
    10 = @{ D = 22.119;  F = 1.111 }
    11 = @{ D = 22.119;  F = 0.804 }
    12 = @{ D = 25.525;  F = 0.683 }
    13 = @{ D = 25.525;  F = 1.273 }
    14 = @{ D = 25.525;  F = 1.175 }
    15 = @{ D = 25.525;  F = 0.869 }
    16 = @{ D = 23.338;  F = 1.202 }
    17 = @{ D = 23.338;  F = 0.672 }
    18 = @{ D = 23.338;  F = 1.315 }
    19 = @{ D = 23.338;  F = 0.811 }
    20 = @{ D = 12.594;  F = 1.295 }
    21 = @{ D = 12.594;  F = 0.972 }
    22 = @{ D = 12.594;  F = 1.235 }
    23 = @{ D = 12.594;  F = 0.862 }
    24 = @{ D = 12.594;  F = 0.637 }
    25 = @{ D = 11.832;  F = 0.912 }
    26 = @{ D = 11.832;  F = 1.496 }
    27 = @{ D = 11.832;  F = 0.812 }
    28 = @{ D = 11.832;  F = 0.78 }
    29 = @{ D = 7.066;   F = 1.555 }
    30 = @{ D = 7.066;   F = 0.898 }
    31 = @{ D = 7.066;   F = 0.931 }
    32 = @{ D = 7.066;   F = 0.616 }
    33 = @{ D = 7.189;   F = 0.87 }
    34 = @{ D = 7.189;   F = 1.246 }
    35 = @{ D = 7.189;   F = 1.242 }
    36 = @{ D = 7.189;   F = 1.012 }
    37 = @{ D = 7.189;   F = 0.629 }
    38 = @{ D = 12.776;  F = 0.384 }
    39 = @{ D = 12.776;  F = 0.523 }
    40 = @{ D = 12.776;  F = 1.382 }
    41 = @{ D = 12.776;  F = 1.71 }
    42 = @{ D = 21.13;   F = 0.53 }
    43 = @{ D = 21.13;   F = 0.888 }
    44 = @{ D = 21.13;   F = 1.489 }
    45 = @{ D = 21.13;   F = 1.093 }
    46 = @{ D = 30.238;  F = 0.894 }
    47 = @{ D = 30.238;  F = 1 }
    48 = @{ D = 30.238;  F = 0.927 }
    49 = @{ D = 30.238;  F = 1.18 }
    50 = @{ D = 31.39;   F = 0.903 }
    51 = @{ D = 31.39;   F = 1.243 }
    52 = @{ D = 31.39;   F = 1.002 }
    53 = @{ D = 31.39;   F = 0.852 }
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row].D
    $ws.Cells.Item($row, 6).Value = $updates[$row].F
}
